# Mask the actual student names on the "OMIS 482" sheet with single
# letters (A-Y) so individual students can no longer be identified.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OMIS 482")

$letters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")

for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $letters[$i]
}

[void]$ws.Range("A27").Select()
